# Update tier_score_sheet (raw performance / base_tier_score values)
$wb = $excel.ActiveWorkbook

$wsTier = $wb.Worksheets.Item("tier_score_sheet")

$wsTier.Range("B3").Value = 400
$wsTier.Range("C3").Value = 4

$wsTier.Range("B4").Value = 500
$wsTier.Range("C4").Value = 5

$wsTier.Range("B5").Value = 500
$wsTier.Range("C5").Value = 5

$wsTier.Range("B6").Value = 550
$wsTier.Range("C6").Value = 5.5

$wsTier.Range("C8").Value = 7.5

$wsTier.Range("C9").Value = 8

# Update overall_tier_scores (performance / base_tier_score + downstream
# calculated columns recomputed to match the new inputs)
$wsOverall = $wb.Worksheets.Item("overall_tier_scores")

# Row 3 (Intel Arc A770 8GB)
$wsOverall.Range("C3").Value = 400
$wsOverall.Range("D3").Value = 4
$wsOverall.Range("M3").Value = 4
$wsOverall.Range("O3").Value = -0.2
$wsOverall.Range("P3").Value = 3.8

# Row 4 (Intel Arc A770 16GB)
$wsOverall.Range("C4").Value = 500
$wsOverall.Range("D4").Value = 5
$wsOverall.Range("L4").Value = 1
$wsOverall.Range("M4").Value = 6
$wsOverall.Range("O4").Value = 0.7500000000000001
$wsOverall.Range("P4").Value = 5.75

# Row 5 (Geforce RTX 3060 8GB)
$wsOverall.Range("C5").Value = 500
$wsOverall.Range("D5").Value = 5
$wsOverall.Range("L5").Value = 0.5
$wsOverall.Range("M5").Value = 5.5
$wsOverall.Range("P5").Value = 5

# Row 6 (Geforce RTX 3060 12GB)
$wsOverall.Range("C6").Value = 550
$wsOverall.Range("D6").Value = 5.5
$wsOverall.Range("L6").Value = 1.1
$wsOverall.Range("M6").Value = 6.6
$wsOverall.Range("O6").Value = 0.55
$wsOverall.Range("P6").Value = 6.05

# Row 8 (Geforce RTX 3080 10GB)
$wsOverall.Range("D8").Value = 7.5
$wsOverall.Range("L8").Value = 1.125
$wsOverall.Range("M8").Value = 8.625
$wsOverall.Range("O8").Value = 0.375
$wsOverall.Range("P8").Value = 7.875

# Row 9 (Geforce RTX 3080 12GB)
$wsOverall.Range("D9").Value = 8
$wsOverall.Range("L9").Value = 1.6
$wsOverall.Range("M9").Value = 9.6
$wsOverall.Range("O9").Value = 0.8
$wsOverall.Range("P9").Value = 8.800000000000001
